# Auto-generated edit script applying numeric cell updates described by the diff
# for Sheets/Jenova_Profits.xlsx (workbook with per-profession sheets).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 53832.58
$ws.Range("I28").Value = 53832.58
$ws.Range("K28").Value = 53832.58
$ws.Range("M28").Value = -53347.58
$ws.Range("H39").Value = 477.14285
$ws.Range("I39").Value = 198
$ws.Range("J39").Value = 632.2222
$ws.Range("K39").Value = 594
$ws.Range("L39").Value = 1896.6666
$ws.Range("M39").Value = -298
$ws.Range("N39").Value = -2488.6666
$ws.Range("H62").Value = 13895290
$ws.Range("I62").Value = 41670400
$ws.Range("J62").Value = 7733.8335
$ws.Range("K62").Value = 41670400
$ws.Range("L62").Value = 7733.8335
$ws.Range("M62").Value = -41669776
$ws.Range("N62").Value = -8981.833500000001
$ws.Range("H64").Value = 6844.75
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752
$ws.Range("H65").Value = 13895290
$ws.Range("I65").Value = 41670400
$ws.Range("J65").Value = 7733.8335
$ws.Range("K65").Value = 208352000
$ws.Range("L65").Value = 38669.1675
$ws.Range("M65").Value = -208348880
$ws.Range("N65").Value = -44909.1675
$ws.Range("H67").Value = 6844.75
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142
$ws.Range("H76").Value = 83339810
$ws.Range("I76").Value = 6203.5
$ws.Range("K76").Value = 6203.5
$ws.Range("M76").Value = -5888.5
$ws.Range("H79").Value = 83339810
$ws.Range("I79").Value = 6203.5
$ws.Range("K79").Value = 6203.5
$ws.Range("M79").Value = -5111.5
$ws.Range("H98").Value = 4249.875
$ws.Range("I98").Value = 3999.8333
$ws.Range("K98").Value = 3999.8333
$ws.Range("M98").Value = -2501.8333
$ws.Range("H122").Value = 4249.875
$ws.Range("I122").Value = 3999.8333
$ws.Range("K122").Value = 11999.4999
$ws.Range("M122").Value = -9549.499899999999
$ws.Range("H132").Value = 13348.028
$ws.Range("I132").Value = 2099.4062
$ws.Range("J132").Value = 133333.33
$ws.Range("K132").Value = 6298.2186
$ws.Range("L132").Value = 399999.99
$ws.Range("M132").Value = -3768.2186
$ws.Range("N132").Value = -405059.99
$ws.Range("H135").Value = 627342.4
$ws.Range("I135").Value = 771076.9
$ws.Range("J135").Value = 4493
$ws.Range("K135").Value = 6939692.100000001
$ws.Range("L135").Value = 40437
$ws.Range("M135").Value = -6937157.100000001
$ws.Range("N135").Value = -45507
$ws.Range("H137").Value = 3782
$ws.Range("I137").Value = 4055.6072
$ws.Range("K137").Value = 12166.8216
$ws.Range("M137").Value = -9616.821599999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1785.1177
$ws.Range("I74").Value = 1695.5834
$ws.Range("K74").Value = 1695.5834
$ws.Range("M74").Value = -821.5834
$ws.Range("H77").Value = 1785.1177
$ws.Range("I77").Value = 1695.5834
$ws.Range("K77").Value = 8477.916999999999
$ws.Range("M77").Value = -4109.916999999999
$ws.Range("H102").Value = 1808.8572
$ws.Range("I102").Value = 1808.8572
$ws.Range("K102").Value = 1808.8572
$ws.Range("M102").Value = -186.8571999999999
$ws.Range("H132").Value = 3779.6938
$ws.Range("I132").Value = 3490.6191
$ws.Range("J132").Value = 5514.143
$ws.Range("K132").Value = 10471.8573
$ws.Range("L132").Value = 16542.429
$ws.Range("M132").Value = -7941.8573
$ws.Range("N132").Value = -21602.429

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29414236
$ws.Range("I20").Value = 55558060
$ws.Range("K20").Value = 55558060
$ws.Range("M20").Value = -55557813
$ws.Range("H134").Value = 38045.484
$ws.Range("I134").Value = 3456.2964
$ws.Range("K134").Value = 10368.8892
$ws.Range("M134").Value = -7833.889200000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 46811.78
$ws.Range("I31").Value = 1111.1818
$ws.Range("J31").Value = 88704
$ws.Range("K31").Value = 1111.1818
$ws.Range("L31").Value = 88704
$ws.Range("M31").Value = -816.1818000000001
$ws.Range("N31").Value = -89294
$ws.Range("H34").Value = 46811.78
$ws.Range("I34").Value = 1111.1818
$ws.Range("J34").Value = 88704
$ws.Range("K34").Value = 1111.1818
$ws.Range("L34").Value = 88704
$ws.Range("M34").Value = -909.1818000000001
$ws.Range("N34").Value = -89108
$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("M75").Value = -36996
$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("M78").Value = -114984
$ws.Range("H99").Value = 4363.773
$ws.Range("I99").Value = 3134.0833
$ws.Range("J99").Value = 5839.4
$ws.Range("K99").Value = 3134.0833
$ws.Range("L99").Value = 5839.4
$ws.Range("M99").Value = -1636.0833
$ws.Range("N99").Value = -8835.4
$ws.Range("H105").Value = 1010.25
$ws.Range("I105").Value = 1010
$ws.Range("J105").Value = 1011
$ws.Range("K105").Value = 1010
$ws.Range("L105").Value = 1011
$ws.Range("M105").Value = 737
$ws.Range("N105").Value = -4505
$ws.Range("H107").Value = 485.7857
$ws.Range("I107").Value = 446.69232
$ws.Range("K107").Value = 446.69232
$ws.Range("M107").Value = 1473.30768
$ws.Range("H122").Value = 4306.75
$ws.Range("I122").Value = 4530
$ws.Range("J122").Value = 4232.3335
$ws.Range("K122").Value = 13590
$ws.Range("L122").Value = 12697.0005
$ws.Range("M122").Value = -11140
$ws.Range("N122").Value = -17597.0005
$ws.Range("H126").Value = 4363.773
$ws.Range("I126").Value = 3134.0833
$ws.Range("J126").Value = 5839.4
$ws.Range("K126").Value = 9402.249899999999
$ws.Range("L126").Value = 17518.2
$ws.Range("M126").Value = -6932.249899999999
$ws.Range("N126").Value = -22458.2
$ws.Range("H134").Value = 258457.52
$ws.Range("I134").Value = 2101.1316
$ws.Range("K134").Value = 6303.3948
$ws.Range("M134").Value = -3768.3948

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 326377.3
$ws.Range("J4").Value = 1000000
$ws.Range("L4").Value = 3000000
$ws.Range("N4").Value = -3000224
$ws.Range("H56").Value = 5979.6
$ws.Range("I56").Value = 5979.6
$ws.Range("K56").Value = 5979.6
$ws.Range("M56").Value = -5449.6
$ws.Range("H97").Value = 732.75
$ws.Range("I97").Value = 953.6667
$ws.Range("K97").Value = 2861.0001
$ws.Range("M97").Value = -2365.0001
$ws.Range("H122").Value = 39671.5
$ws.Range("I122").Value = 809.5
$ws.Range("J122").Value = 63960.25
$ws.Range("K122").Value = 7285.5
$ws.Range("L122").Value = 575642.25
$ws.Range("M122").Value = -4835.5
$ws.Range("N122").Value = -580542.25
$ws.Range("H132").Value = 504286.3
$ws.Range("I132").Value = 112128.89
$ws.Range("J132").Value = 775779.9399999999
$ws.Range("K132").Value = 1009160.01
$ws.Range("L132").Value = 6982019.459999999
$ws.Range("M132").Value = -1006630.01
$ws.Range("N132").Value = -6987079.459999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8350000
$ws.Range("I11").Value = 7525000
$ws.Range("K11").Value = 7525000
$ws.Range("M11").Value = -7524861
$ws.Range("H14").Value = 63742424
$ws.Range("I14").Value = 71782824
$ws.Range("K14").Value = 71782824
$ws.Range("M14").Value = -71782656
$ws.Range("H52").Value = 30031.334
$ws.Range("J52").Value = 30031.334
$ws.Range("L52").Value = 30031.334
$ws.Range("N52").Value = -30549.334
$ws.Range("H97").Value = 979.4737
$ws.Range("I97").Value = 998.86664
$ws.Range("K97").Value = 998.86664
$ws.Range("M97").Value = -502.86664
$ws.Range("H126").Value = 3275.7693
$ws.Range("I126").Value = 2926.1667
$ws.Range("J126").Value = 3575.4285
$ws.Range("K126").Value = 8778.500100000001
$ws.Range("L126").Value = 10726.2855
$ws.Range("M126").Value = -6308.500100000001
$ws.Range("N126").Value = -15666.2855
$ws.Range("H132").Value = 61618.473
$ws.Range("I132").Value = 8573.588
$ws.Range("K132").Value = 25720.764
$ws.Range("M132").Value = -23190.764

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15083.5
$ws.Range("I7").Value = 22268.666
$ws.Range("K7").Value = 22268.666
$ws.Range("M7").Value = -22156.666
$ws.Range("H16").Value = 201
$ws.Range("I16").Value = 201
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 201
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -31
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 5468.364
$ws.Range("I40").Value = 3367
$ws.Range("J40").Value = 7990
$ws.Range("K40").Value = 3367
$ws.Range("L40").Value = 7990
$ws.Range("M40").Value = -3231
$ws.Range("N40").Value = -8262
$ws.Range("H122").Value = 4536.6875
$ws.Range("I122").Value = 3932.7778
$ws.Range("K122").Value = 11798.3334
$ws.Range("M122").Value = -9348.3334
$ws.Range("H126").Value = 15083.5
$ws.Range("I126").Value = 22268.666
$ws.Range("K126").Value = 66805.99800000001
$ws.Range("M126").Value = -64335.99800000001
$ws.Range("H136").Value = 391245.47
$ws.Range("I136").Value = 629111.4399999999
$ws.Range("K136").Value = 1887334.32
$ws.Range("M136").Value = -1884784.32

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 719.5925999999999
$ws.Range("I107").Value = 762.9545000000001
$ws.Range("K107").Value = 2288.8635
$ws.Range("M107").Value = -368.8635000000004
$ws.Range("H132").Value = 51651.047
$ws.Range("I132").Value = 3073.25
$ws.Range("K132").Value = 9219.75
$ws.Range("M132").Value = -6689.75
$ws.Range("H136").Value = 13416736
$ws.Range("I136").Value = 16670220
$ws.Range("K136").Value = 50010660
$ws.Range("M136").Value = -50008110

Write-Output "Applied 246 value updates and 1 cell clear across 8 sheets."